# Update countries & provincias Spain
# This script applies the "data refresh" edit captured by the diff:
#  - the "last updated" timestamp in A1 changes
#  - several countries' case counts (columns B-H) are refreshed
#  - because the sheet is sorted by total cases (col B) descending, a
#    handful of adjacent country rows swap rank / labels as a result

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Cells.Item(1,1).Value = "Datos actualizados a 3 de Julio de 2020 a las 15:58"

# --- Estados Unidos (row 4) --------------------------------------------
$ws.Cells.Item(4,2).Value = 2837726
$ws.Cells.Item(4,3).Value = 537
$ws.Cells.Item(4,5).Value = 1514384
$ws.Cells.Item(4,7).Value = 19
$ws.Cells.Item(4,8).Value = 131504

# --- India (row 7) ------------------------------------------------------
$ws.Cells.Item(7,2).Value = 632998
$ws.Cells.Item(7,3).Value = 5830
$ws.Cells.Item(7,4).Value = 383516
$ws.Cells.Item(7,5).Value = 231171
$ws.Cells.Item(7,7).Value = 86
$ws.Cells.Item(7,8).Value = 18311

# --- Arabia Saudita (row 17) --------------------------------------------
$ws.Cells.Item(17,2).Value = 201801
$ws.Cells.Item(17,3).Value = 4193
$ws.Cells.Item(17,4).Value = 140614
$ws.Cells.Item(17,5).Value = 59385
$ws.Cells.Item(17,7).Value = 50
$ws.Cells.Item(17,8).Value = 1802

# --- Argentina (row 28) --------------------------------------------------
$ws.Cells.Item(28,4).Value = 25224
$ws.Cells.Item(28,5).Value = 43314
$ws.Cells.Item(28,7).Value = 18
$ws.Cells.Item(28,8).Value = 1403

# --- Irak (row 33) --------------------------------------------------------
$ws.Cells.Item(33,2).Value = 56020
$ws.Cells.Item(33,3).Value = 2312
$ws.Cells.Item(33,4).Value = 29600
$ws.Cells.Item(33,5).Value = 24158
$ws.Cells.Item(33,7).Value = 102
$ws.Cells.Item(33,8).Value = 2262

# --- Japon / Azerbaiyan swap (rows 57-58) ---------------------------------
$ws.Cells.Item(57,1).Value = "Azerbaiyan"
$ws.Cells.Item(57,2).Value = 19267
$ws.Cells.Item(57,3).Value = 583
$ws.Cells.Item(57,4).Value = 10820
$ws.Cells.Item(57,5).Value = 8212
$ws.Cells.Item(57,7).Value = 7
$ws.Cells.Item(57,8).Value = 235

$ws.Cells.Item(58,1).Value = "Japon"
$ws.Cells.Item(58,2).Value = 18874
$ws.Cells.Item(58,4).Value = 16772
$ws.Cells.Item(58,5).Value = 1127
$ws.Cells.Item(58,8).Value = 975

# --- Serbia / Nepal swap (rows 62-63) -------------------------------------
$ws.Cells.Item(62,1).Value = "Serbia"
$ws.Cells.Item(62,2).Value = 15504
$ws.Cells.Item(62,3).Value = 309
$ws.Cells.Item(62,4).Value = 13064
$ws.Cells.Item(62,5).Value = 2142
$ws.Cells.Item(62,7).Value = 11
$ws.Cells.Item(62,8).Value = 298

$ws.Cells.Item(63,1).Value = "Nepal"
$ws.Cells.Item(63,2).Value = 15259
$ws.Cells.Item(63,3).Value = 740
$ws.Cells.Item(63,4).Value = 6143
$ws.Cells.Item(63,5).Value = 9084
$ws.Cells.Item(63,7).Value = 1
$ws.Cells.Item(63,8).Value = 32

# --- Finlandia (row 78) ---------------------------------------------------
$ws.Cells.Item(78,5).Value = 213
$ws.Cells.Item(78,7).Value = 1
$ws.Cells.Item(78,8).Value = 329

# --- Kenia (row 79) --------------------------------------------------------
$ws.Cells.Item(79,5).Value = 4925
$ws.Cells.Item(79,7).Value = 2
$ws.Cells.Item(79,8).Value = 154

# --- Libano / Lituania swap (rows 116-117) --------------------------------
$ws.Cells.Item(116,1).Value = "Libano"
$ws.Cells.Item(116,2).Value = 1830
$ws.Cells.Item(116,3).Value = 34
$ws.Cells.Item(116,4).Value = 1292
$ws.Cells.Item(116,5).Value = 503
$ws.Cells.Item(116,7).Value = 0
$ws.Cells.Item(116,8).Value = 35

$ws.Cells.Item(117,1).Value = "Lituania"
$ws.Cells.Item(117,2).Value = 1828
$ws.Cells.Item(117,3).Value = 3
$ws.Cells.Item(117,4).Value = 1539
$ws.Cells.Item(117,5).Value = 210
$ws.Cells.Item(117,7).Value = 1
$ws.Cells.Item(117,8).Value = 79

# --- Hong Kong (row 127) ---------------------------------------------------
$ws.Cells.Item(127,2).Value = 1248
$ws.Cells.Item(127,3).Value = 5
$ws.Cells.Item(127,4).Value = 1125

# --- Reunion (row 155) ------------------------------------------------------
$ws.Cells.Item(155,2).Value = 533
$ws.Cells.Item(155,3).Value = 2
$ws.Cells.Item(155,5).Value = 59

# --- Siria / Angola swap (rows 162-163) -------------------------------------
$ws.Cells.Item(162,1).Value = "Siria"
$ws.Cells.Item(162,2).Value = 328
$ws.Cells.Item(162,3).Value = 16
$ws.Cells.Item(162,4).Value = 123
$ws.Cells.Item(162,5).Value = 195
$ws.Cells.Item(162,7).Value = 1
$ws.Cells.Item(162,8).Value = 10

$ws.Cells.Item(163,1).Value = "Angola"
$ws.Cells.Item(163,2).Value = 315
$ws.Cells.Item(163,4).Value = 97
$ws.Cells.Item(163,5).Value = 201
$ws.Cells.Item(163,8).Value = 17

# --- Fiyi / Dominica tie-break swap (rows 205-206), values unchanged -------
$ws.Cells.Item(205,1).Value = "Fiyi"
$ws.Cells.Item(206,1).Value = "Dominica"

# --- Groenlandia / Islas Malvinas tie-break swap (rows 209-210) ------------
$ws.Cells.Item(209,1).Value = "Groenlandia"
$ws.Cells.Item(210,1).Value = "Islas Malvinas"
